$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update C6 (the "1.4" row description): extend the bold lead-in run's text ---
$c6 = $ws.Range("C6")

$oldBoldLen = 38
$newBoldText = "Officieel klaar met de user stories!! in versie 2.1 zijn echter alle bugs er pas uit.  "
$newBoldLen = $newBoldText.Length
$secondRunLen = 434

$c6.Characters(1, $oldBoldLen).Text = $newBoldText

# Re-assert the run-level bold formatting so the rich text stays split into
# the same two runs (bold lead-in + regular continuation) after the edit.
$c6.Characters(1, $newBoldLen).Font.Bold = $true
$c6.Characters(($newBoldLen + 1), $secondRunLen).Font.Bold = $false

# --- Add new row 8 entry (hour 7, version 2.1) ---
$ws.Range("A8").Value = 7
# Force text so "2.1" isn't auto-coerced into the number 2.1 (leading
# apostrophe marks it as text, matching how the other "x.y" version
# cells in column B are stored as shared strings).
$ws.Range("B8").Value = "'2.1"
$ws.Range("C8").Value = "Ik heb 2 extra opties toegevoegd om performance te verbeteren bij grote mazes. Eentje op het licht uit te zetten, en de andere om de particles uit te zetten. De  `"Pretty mode`" functie heet nu `"Fast mode`" voor verduiderlijking. Hiernaast heb ik code geschreven zodat de oude cells die nog gegenereerd werden allemaal verwijderd worden, doormiddel van StopCoroutine. Hiervoor  gebeurde dit niet altijd."

# Match the styling used by the other data rows (B: s=2, C: s=3) by
# pasting formats from the row above — this also clears the transient
# "quote prefix" style that entering "'2.1" applied to B8.
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# --- Update the active selection to C13 ---
$ws.Range("C13").Select()
